$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor values for rows 2-21, columns C-H (ax, ay, az, gx, gy, gz)
# col order: C=3, D=4, E=5, F=6, G=7, H=8
$data = @{
    2 = @(-1.487755507230762; -10.35288465023041; 3.703868120908739; -2.454914808273315; -0.2842050492763519; -2.869386672973633)
    3 = @(-4.7107778787613; -8.213166236877427; 5.28016042709351; 2.025781869888306; 1.884366512298584; 0.9054549932479858)
    4 = @(-6.86208564043045; -5.234421730041505; 4.114273145794865; 1.317789077758789; -1.244027018547058; 1.031140565872192)
    5 = @(-6.348576307296753; -6.586597442626955; 1.944473057985312; 0.7573201656341553; -0.4138612151145935; 0.1753183305263519)
    6 = @(-4.850189305841925; -7.467463970184327; 6.846671938896179; 0.5755877494812012; -0.2176207453012466; -1.24479067325592)
    7 = @(0.7637936472892619; -3.36792755126954; 6.735066413879379; 0.6409503817558289; -4.134488582611084; 0.0209221355617046)
    8 = @(-5.654898561537283; -7.294652104377757; -0.4670590609312262; -1.424232363700867; -3.165656328201294; -0.2144137024879455)
    9 = @(1.988154411315929; -8.051345348358154; -4.82868674397469; -2.829985857009888; -4.054770946502686; 1.877952456474304)
    10 = @(1.445369705557815; -8.260547280311584; -3.376170933246607; 4.016286373138428; -3.601814270019531; -2.621833562850952)
    11 = @(-2.370766460895542; -7.741124629974364; -6.762903213500985; -0.9879218339920044; 4.384180068969727; 1.598176240921021)
    12 = @(-2.55821630358695; -8.715493917465219; -6.844864040613141; -3.011718273162842; 2.93963623046875; 0.5337435007095337)
    13 = @(-0.8884068131446955; -10.84955549240112; 1.661297619342809; -0.4014911949634552; 0.9065240025520324; -0.6395758986473083)
    14 = @(-3.758049488067625; -8.76593214273451; 3.196124792099005; -0.99250328540802; 2.58121109008789; -0.3058907687664032)
    15 = @(-3.226068019866951; -4.671105861663801; 4.809621334075939; -2.1601722240448; 2.302809238433838; -0.3118467032909393)
    16 = @(-4.81333899497986; -0.7696201205253579; 7.391847074031832; -0.5739079117774963; -1.005178809165955; 0.2600758671760559)
    17 = @(-0.6115292310714668; -3.946480035781865; 7.73033595085144; -1.500132322311401; 1.276097536087036; -1.762957096099854)
    18 = @(-4.580467939376843; -2.483635365962976; 9.400022864341736; -1.534646153450012; -1.36818540096283; 0.8046622276306152)
    19 = @(-3.186845898628217; 3.314929008483903; 7.02214622497557; -1.083674907684326; 1.494787216186523; 0.7244861721992493)
    20 = @(-0.8589091598987735; 5.796535491943349; 4.699775040149696; 1.442405581474304; -1.076344609260559; 0.901637077331543)
    21 = @(-6.590956926345826; -0.2845994234085083; 5.918066263198853; 4.085466861724854; -3.372739791870117; 4.743368625640869)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 3 + $i).Value2 = $vals[$i]
    }
}

# Remove the now-extra last row (old row 22) and shrink the sheet dimension
$ws.Rows.Item(22).Delete()

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
